$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("best-feasible-slns")

# Updated "actual" values in column C (rows 2-31, row 7 unchanged)
$updates = @{
    2  = 23237
    3  = 21852
    4  = 22034
    5  = 23393
    6  = 23137
    8  = 24969
    9  = 23043
    10 = 22947
    11 = 22728
    12 = 41982
    13 = 39655
    14 = 39643
    15 = 42029
    16 = 40035
    17 = 39455
    18 = 40764
    19 = 42271
    20 = 41819
    21 = 43888
    22 = 59464
    23 = 62127
    24 = 59470
    25 = 60226
    26 = 60007
    27 = 58747
    28 = 61728
    29 = 61727
    30 = 59348
    31 = 60447
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}

# Update the active-cell selection recorded in the sheet view
$ws.Range("H9").Select()
